# Update "想去人数" (interested-count) figures for two events that are
# tracked on both the "展览" sheet and the aggregated "全部类型" sheet.
#   F3: 109 -> 110
#   F6: 13  -> 14

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 110
    $ws.Range("F6").Value = 14
}
